# Update the "想去人数" (F column) figures on both the "展览" sheet and the
# "全部类型" sheet. Both sheets share the same rows/values for this column,
# and both need to be bumped identically.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 1175
    4  = 41
    6  = 170
    10 = 5465
    11 = 4872
    16 = 198
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
